$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing values
$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = 3

# Remove row 4 entirely (A4=3, B4=1) - shift cells up, deleting the row
$ws.Rows.Item(4).Delete()
